# Import Master data dan Add Sensor pemasangan dan maintennace
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update IMEI numbers for the two existing rows
$ws.Range("C2").Value = 873652413245321
$ws.Range("C3").Value = 152635423178765

# Clear the "Status ownership" value for row 2 (was "-")
$ws.Range("G2").ClearContents()

# Move the active selection to D5
$ws.Range("D5").Select()
